$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "3376 Royal Freemasons Coppin Centre Melbourne",
    "3398 BlueCross Elly Kay Mordialloc",
    "3564 Waverley Valley Aged Care Glen Waverley",
    "3601 Baptcare Westhaven community",
    "3647 Aurrum Aged Care Reservoir",
    "3653 Fronditha Thalpori St Albans Aged Care",
    "3975 Aurrum Aged Care Brunswick West",
    "4257 BlueCross The Gables Camberwell",
    "4295 Hope Aged Care Sunshine West",
    "4314 Estia Health Ardeer",
    "44095 Myrniong Primary School Myrniong",
    "44304 Brighton Primary School Brighton",
    "44404 Castlemaine North Primary School Castlemaine",
    "44490 Armadale Primary School Armadale",
    "44593 Torquay P-6 College Torquay",
    "44620 Canterbury Primary School Canterbury",
    "44623 Brunswick North Primary School Brunswick West",
    "44745 Briar Hill Primary School Briar Hill",
    "44765 Strathmore Primary School Strathmore",
    "44799 Eastwood Primary School Ringwood East",
    "44960 Thomastown West Primary School",
    "45013 Gladstone Views Primary School",
    "45147 Maramba Primary School Narre Warren",
    "45168 Ranfurly Primary School Mildura",
    "45181 Courtenay Gardens Primary School Cranbourne North",
    "45257 Roxburgh Rise Primary School Roxburgh Park",
    "45305 Lockington Consolidated School Lockington",
    "45719 St Joseph's Primary School Numurkah",
    "4574 Village Glen Aged Care Residences Mornington",
    "45757 Saint Joseph's Primary School Warragul",
    "45764 Our Lady Help of Christian's Primary School Brunswick East",
    "45858 St Bernard's Primary Coburg",
    "45861 St Oliver Plunkett Primary School Pascoe Vale",
    "45958 Ave Maria College Aberfeldie Workplace",
    "45975 St Thomas More Primary School Hadfield",
    "46074 St Justin's Catholic Primary School Wheelers Hill",
    "46078 Corpus Christi Primary School Werribee",
    "46086 St Kevin's Primary School Hampton Park",
    "46104 Clairvaux Catholic School Belmont",
    "46135 Wesley College Junior School St Kilda Road Melbourne",
    "46208 Mount Scopus Memorial College Gandel Campus Burwood",
    "46327 Victory Christian College Strathdale",
    "50279 Dallas Brooks Community Primary School Dallas",
    "51529 Sirius College Primary School Dallas",
    "52390 Our Lady of the Way Catholic Primary School Wallan",
    "52694 Pakenham Primary School Pakenham",
    "Australian Radio Network Richmond",
    "Ballarat Freedom Protest",
    "Brandon Park Primary School Wheelers Hill",
    "Cardinia Waters Retirement Village Pakenham",
    "Confirmed Omicron Sircuit Bar Fitzroy",
    "Confirmed Omicron Variant The Peel Hotel Collingwood",
    "Goodstart Early Learning Preston",
    "Greendale Hotel Greendale",
    "JBS Australia Brooklyn",
    "Kororoit Creek Primary School Burnside Heights",
    "Oakleigh South Primary School Oakleigh South",
    "PGL Camp Rumbug Foster North",
    "Rosebud Primary School Rosebud",
    "Social Gathering 11 Dec Windsor",
    "St Brigid's Parish Primary School Mordialloc",
    "St Clare's Primary School Officer",
    "St Mary's Parish Primary School",
    "St Vincents Hospital Melbourne Emergency Department Fitzroy",
    "StarTrack Tullamarine",
    "The George Lounge St Kilda",
    "Thomastown West Primary School Camp Doxa's Malmsbury",
    "V & G construction site San Lorenzo Wine & Dining"
)

$values = @(
    10,
    30,
    13,
    13,
    12,
    22,
    13,
    16,
    16,
    17,
    13,
    12,
    49,
    27,
    35,
    18,
    28,
    22,
    10,
    37,
    14,
    37,
    13,
    25,
    11,
    10,
    33,
    12,
    10,
    12,
    12,
    27,
    10,
    22,
    11,
    15,
    34,
    11,
    10,
    10,
    12,
    14,
    13,
    10,
    33,
    22,
    14,
    10,
    10,
    14,
    15,
    17,
    12,
    21,
    29,
    24,
    11,
    47,
    19,
    22,
    12,
    10,
    11,
    11,
    16,
    16,
    20,
    12
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Rows.Item(70).Delete()

